$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "IWChecker"
$ws.Range("B2").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("C2").Value = "19/05/2021 06:31:27"

$ws.Range("A3").Value = "IWMaker"
$ws.Range("B3").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("C3").Value = "19/05/2021 06:30:55"

$ws.Range("A4").Value = "OCMChecker"
$ws.Range("B4").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("C4").Value = "02/06/2021 15:55:37"

$ws.Range("A5").Value = "OCMMaker"
$ws.Range("B5").Value = "EC2AMAZ-N8SAHHO\Administrator"
$ws.Range("C5").Value = "03/06/2021 15:13:14"

$wb.Save()
